$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 1)
# Two new "bold" aggregate label columns are inserted after the blank
# A1 cell, and the country columns are reordered (Poland/Spain/UK/
# Switzerland moved next to the other European countries, and
# "United States" is renamed "USA"). Two new country columns
# (Saudi Arabia, USA) are appended at the end (L, M).
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "`$ bold(All)"
$ws.Range("C1").Value = "`$ bold(Europe)"
$ws.Range("D1").Value = "France"
$ws.Range("E1").Value = "Germany"
$ws.Range("F1").Value = "Italy"
$ws.Range("G1").Value = "Poland"
$ws.Range("H1").Value = "Spain"
$ws.Range("I1").Value = "United Kingdom"
$ws.Range("J1").Value = "Switzerland"
$ws.Range("K1").Value = "Japan"
$ws.Range("L1").Value = "Saudi Arabia"
$ws.Range("M1").Value = "USA"

# ---------------------------------------------------------------------
# Row labels (column A, rows 2-5) - the `*High/*Mid/*Low` footnote
# labels gain markup wrapping (bold()/paste()) and a couple of stray
# trailing parentheses are dropped.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "`$ paste('Supports the GCS if its other members* cover 64-72%`nof world emissions*', bold(High), ' : Global South + China +`nEU + various HICs (UK, Japan, South Korea, Canada...)')"
$ws.Range("A3").Value = "`$ paste('Supports the GCS if its other members* cover 64-72% of world emissions,`n*', bold('High color'), ' : High + Distributive effects displayed using colors on world map')"
$ws.Range("A4").Value = "Supports the GCS if its other members* cover 56% of world emissions`n*Mid: Global South + China"
$ws.Range("A5").Value = "Supports the GCS if its other members* cover 25-33% of world emissions`n*Low: Global South + EU"

# ---------------------------------------------------------------------
# Data values (rows 2-5, columns B-M)
# ---------------------------------------------------------------------
$data = @(
    @(0.684797846013269, 0.730390283144958, 0.721888822593033, 0.793220476440396, 0.806061731494686, 0.725315089483971, 0.709725010591962, 0.60367746807593, 0.802921625022916, 0.703649688948398, 0.892630395702026, 0.613190097859686),
    @(0.631922747338702, 0.68464172802944, 0.60414441640184, 0.692000384560851, 0.764023624110181, 0.593338755591008, 0.764978780580193, 0.687158186636906, 0.678231464642976, 0.630666968689808, 0.868668628355021, 0.557762638731718),
    @(0.679004869341497, 0.74289069142581, 0.771098600012077, 0.769236558846297, 0.83596654741558, 0.566504911211525, 0.767597594944841, 0.670603338366483, 0.633953440855902, 0.682575327407304, 0.876160911546928, 0.570538630080154),
    @(0.645460623471237, 0.681993956456158, 0.762733353729202, 0.622594602050181, 0.820338845276832, 0.618966284979516, 0.690022256508678, 0.612454933195732, 0.631254517516158, 0.670645142503371, 0.868214391306789, 0.573626811841157)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
